$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 605
$ws.Cells.Item(4, 6).Value = 6434
$ws.Cells.Item(5, 6).Value = 731
$ws.Cells.Item(6, 6).Value = 1090
$ws.Cells.Item(7, 6).Value = 79
$ws.Cells.Item(8, 6).Value = 502
$ws.Cells.Item(8, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/om0OCpxy1709287210276.jpeg"
$ws.Cells.Item(9, 6).Value = 195
$ws.Cells.Item(10, 6).Value = 26
$ws.Cells.Item(11, 6).Value = 714
$ws.Cells.Item(13, 6).Value = 3
$ws.Cells.Item(14, 6).Value = 85
$ws.Cells.Item(15, 6).Value = 201
$ws.Cells.Item(16, 6).Value = 441
$ws.Cells.Item(17, 6).Value = 50
$ws.Cells.Item(18, 6).Value = 23
$ws.Cells.Item(19, 6).Value = 1419
$ws.Cells.Item(20, 6).Value = 675
$ws.Cells.Item(21, 6).Value = 388
$ws.Cells.Item(22, 6).Value = 403
$ws.Cells.Item(24, 6).Value = 1075
$ws.Cells.Item(25, 6).Value = 155
$ws.Cells.Item(26, 6).Value = 2227
$ws.Cells.Item(28, 6).Value = 107
$ws.Cells.Item(31, 6).Value = 3609
$ws.Cells.Item(33, 6).Value = 637

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 198
$ws.Cells.Item(8, 6).Value = 27
$ws.Cells.Item(9, 6).Value = 712
$ws.Cells.Item(12, 6).Value = 1020
$ws.Cells.Item(14, 6).Value = 111
$ws.Cells.Item(18, 6).Value = 77
$ws.Cells.Item(31, 6).Value = 212
$ws.Cells.Item(34, 6).Value = 49
$ws.Cells.Item(37, 6).Value = 5

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 46
$ws.Cells.Item(4, 6).Value = 1200
$ws.Cells.Item(6, 6).Value = 1581
$ws.Cells.Item(7, 6).Value = 433
$ws.Cells.Item(10, 6).Value = 790

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 1200
$ws.Cells.Item(4, 6).Value = 1581
$ws.Cells.Item(5, 6).Value = 433
$ws.Cells.Item(7, 6).Value = 790
$ws.Cells.Item(8, 6).Value = 605
$ws.Cells.Item(9, 6).Value = 198
$ws.Cells.Item(10, 6).Value = 6434
$ws.Cells.Item(11, 6).Value = 27
$ws.Cells.Item(12, 6).Value = 731
$ws.Cells.Item(13, 6).Value = 1090
$ws.Cells.Item(14, 6).Value = 712
$ws.Cells.Item(15, 6).Value = 79
$ws.Cells.Item(16, 6).Value = 502
$ws.Cells.Item(16, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/om0OCpxy1709287210276.jpeg"
$ws.Cells.Item(17, 6).Value = 195
$ws.Cells.Item(18, 6).Value = 26
$ws.Cells.Item(19, 6).Value = 714
$ws.Cells.Item(21, 6).Value = 111
$ws.Cells.Item(22, 6).Value = 111
$ws.Cells.Item(23, 6).Value = 77
$ws.Cells.Item(25, 6).Value = 85
$ws.Cells.Item(29, 6).Value = 23
$ws.Cells.Item(32, 6).Value = 675
$ws.Cells.Item(33, 6).Value = 388
$ws.Cells.Item(34, 6).Value = 403
$ws.Cells.Item(37, 6).Value = 1075
$ws.Cells.Item(38, 6).Value = 155
$ws.Cells.Item(39, 6).Value = 2227
$ws.Cells.Item(40, 6).Value = 212
$ws.Cells.Item(42, 6).Value = 49
$ws.Cells.Item(45, 6).Value = 107
$ws.Cells.Item(47, 6).Value = 3609
$ws.Cells.Item(48, 6).Value = 5
$ws.Cells.Item(51, 6).Value = 638
